# The diff appends 9 data rows (rows 6-14) of order/pricing data to the
# bottom of "Sheet2" (columns A, C, E, F only - B "NSN Number" and D
# "Delivery Time" are left blank for these rows), growing the sheet's used
# range from A1:F5 to A1:F14.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Each tuple is: PartNumber(A), MinimumOrderQuantity(C), FinalPrice(E), CostPrice(F)
$data = @(
    @(22356568, 1, 112.17,   7.3),
    @(34256354, 1, 200,      4.34),
    @(44637355, 1, 114.2,    74.23),
    @(55261434, 2, 27.67,    18.56),
    @(66464788, 1, 505,      411.14),
    @(77423423, 1, 800,      118.84),
    @(88888856, 1, 32.25,    20.96),
    @(94757647, 1, 812.01,   682.09),
    @(10342423, 1, 28.68,    1.29)
)

$startRow = 6
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]   # A - Part Number
    $ws.Cells.Item($row, 3).Value = $entry[1]   # C - Minimum Order Quantity
    $ws.Cells.Item($row, 5).Value = $entry[2]   # E - Final Price
    $ws.Cells.Item($row, 6).Value = $entry[3]   # F - Cost Price
}
